$d = $word.ActiveDocument

# Locate the paragraph that begins the "Update-" checklist line, e.g.:
#   "Update- Account, Customer, Transaction, Loan, and CreditCard"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("Update- ")) {
        $target = $para
        break
    }
}

$pStart = $target.Range.Start
$pText = $target.Range.Text

# Word offsets (relative to the paragraph) of the "Customer" / "Transaction"
# item names (without their trailing comma).
$custIdx = $pText.IndexOf("Customer,")
$transIdx = $pText.IndexOf("Transaction,")

# The existing "_GoBack" bookmark currently wraps "Account"; it needs to move
# so it wraps "Transaction" instead once that word is marked done. Drop it
# now and re-create it in the right spot below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Mark "Customer" as done (strike-through), matching "Account"'s formatting.
$rCustomer = $d.Range($pStart + $custIdx, $pStart + $custIdx + 8)
$rCustomer.Font.StrikeThrough = $true

# Mark "Transaction" as done (strike-through) too.
$rTransaction = $d.Range($pStart + $transIdx, $pStart + $transIdx + 11)
$rTransaction.Font.StrikeThrough = $true

# Re-insert the "_GoBack" bookmark around "Transaction".
$d.Bookmarks.Add("_GoBack", $rTransaction)
